# Trade #17 closed at 2026-02-17 23:54:13 - unknown UNKNOWN +0.000%
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the new closed trade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1500.33
$wsSummary.Range("B4").Value = 0.33
$wsSummary.Range("B5").Value = 0.39
$wsSummary.Range("B6").Value = 17
$wsSummary.Range("B8").Value = 6
$wsSummary.Range("B9").Value = 58.82

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 100.33
$wsStatus.Range("D6").Value = 17
$wsStatus.Range("E6").Value = 0.33
$wsStatus.Range("F6").Value = 0.33
$wsStatus.Range("G6").Value = 58.82

# ---------------------------------------------------------------------
# New trade row data (shared by "All Trades" and "MarketMaking" sheets)
# ---------------------------------------------------------------------
$newRow = @(17, "2026-02-17", "23:54:07", "MarketMaking", "DOWN", 0.04, 0.02, "CLOSED", -50, -0.02, 100.33, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.18)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 18
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $cell = $ws.Cells.Item($row, $i + 1)
        if ($i -eq 1) {
            # Date column ("2026-02-17") - force text storage instead of
            # letting Excel auto-convert the date-shaped string into a
            # date serial number.
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$i]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newRow[$i]
        }
    }
}
